$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 12).Value = 5331
$ws.Cells.Item(3, 12).Value = 5750
$ws.Cells.Item(4, 12).Value = 1404
$ws.Cells.Item(5, 12).Value = 342
$ws.Cells.Item(6, 12).Value = 4765
$ws.Cells.Item(7, 12).Value = 17592

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(3, 12).Value = 49
$ws.Cells.Item(7, 12).Value = 193

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(2, 12).Value = 338
$ws.Cells.Item(3, 12).Value = 406
$ws.Cells.Item(6, 12).Value = 297
$ws.Cells.Item(7, 12).Value = 1167

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 12).Value = 123
$ws.Cells.Item(3, 12).Value = 156
$ws.Cells.Item(6, 12).Value = 84
$ws.Cells.Item(7, 12).Value = 385

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(4, 12).Value = 54
$ws.Cells.Item(7, 12).Value = 815

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(6, 12).Value = 54
$ws.Cells.Item(7, 12).Value = 246

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(2, 12).Value = 200
$ws.Cells.Item(3, 12).Value = 231
$ws.Cells.Item(6, 12).Value = 183
$ws.Cells.Item(7, 12).Value = 668

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(3, 12).Value = 107
$ws.Cells.Item(7, 12).Value = 341

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Cells.Item(5, 12).Value = 2
$ws.Cells.Item(7, 12).Value = 79

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(7, 12).Value = 573
$ws.Cells.Item(8, 12).Value = 1167
$ws.Cells.Item(9, 12).Value = 101
$ws.Cells.Item(11, 12).Value = 289
$ws.Cells.Item(15, 12).Value = 134
$ws.Cells.Item(18, 12).Value = 122
$ws.Cells.Item(19, 12).Value = 477
$ws.Cells.Item(20, 12).Value = 434
$ws.Cells.Item(24, 12).Value = 45
$ws.Cells.Item(27, 12).Value = 154
$ws.Cells.Item(29, 12).Value = 992
$ws.Cells.Item(30, 12).Value = 79
$ws.Cells.Item(31, 12).Value = 172
$ws.Cells.Item(33, 12).Value = 815
$ws.Cells.Item(36, 12).Value = 223
$ws.Cells.Item(37, 12).Value = 668
$ws.Cells.Item(42, 12).Value = 570
$ws.Cells.Item(43, 12).Value = 130
$ws.Cells.Item(44, 12).Value = 120
$ws.Cells.Item(45, 12).Value = 32
$ws.Cells.Item(51, 12).Value = 223
$ws.Cells.Item(52, 12).Value = 354
$ws.Cells.Item(53, 12).Value = 193
$ws.Cells.Item(57, 12).Value = 62
$ws.Cells.Item(60, 12).Value = 113
$ws.Cells.Item(63, 12).Value = 52
$ws.Cells.Item(65, 12).Value = 341
$ws.Cells.Item(67, 12).Value = 599
$ws.Cells.Item(71, 12).Value = 49
$ws.Cells.Item(76, 12).Value = 275
$ws.Cells.Item(77, 12).Value = 117
$ws.Cells.Item(79, 12).Value = 484
$ws.Cells.Item(83, 12).Value = 385
$ws.Cells.Item(84, 12).Value = 171
$ws.Cells.Item(85, 12).Value = 880
$ws.Cells.Item(88, 12).Value = 192
$ws.Cells.Item(94, 12).Value = 217
$ws.Cells.Item(95, 12).Value = 246
$ws.Cells.Item(96, 12).Value = 199
$ws.Cells.Item(101, 12).Value = 17592

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(3, 12).Value = 45
$ws.Cells.Item(7, 12).Value = 172

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 12).Value = 178
$ws.Cells.Item(7, 12).Value = 599

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(2, 12).Value = 55
$ws.Cells.Item(7, 12).Value = 171

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(2, 12).Value = 299
$ws.Cells.Item(3, 12).Value = 379
$ws.Cells.Item(4, 12).Value = 50
$ws.Cells.Item(6, 12).Value = 247
$ws.Cells.Item(7, 12).Value = 992

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 12).Value = 171
$ws.Cells.Item(7, 12).Value = 477

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(2, 12).Value = 47
$ws.Cells.Item(3, 12).Value = 35
$ws.Cells.Item(7, 12).Value = 120

$ws = $wb.Worksheets.Item("River North")
$ws.Cells.Item(3, 12).Value = 53
$ws.Cells.Item(7, 12).Value = 275

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 12).Value = 193
$ws.Cells.Item(4, 12).Value = 48
$ws.Cells.Item(6, 12).Value = 158
$ws.Cells.Item(7, 12).Value = 570

$ws = $wb.Worksheets.Item("Dunning")
$ws.Cells.Item(2, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 45

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(6, 12).Value = 58
$ws.Cells.Item(7, 12).Value = 199

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(6, 12).Value = 123
$ws.Cells.Item(7, 12).Value = 484

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(2, 12).Value = 136
$ws.Cells.Item(6, 12).Value = 114
$ws.Cells.Item(7, 12).Value = 434

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Cells.Item(4, 12).Value = 11
$ws.Cells.Item(7, 12).Value = 122

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Cells.Item(3, 12).Value = 68
$ws.Cells.Item(7, 12).Value = 223

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 12).Value = 195
$ws.Cells.Item(6, 12).Value = 135
$ws.Cells.Item(7, 12).Value = 573

$ws = $wb.Worksheets.Item("West Loop")
$ws.Cells.Item(3, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 217

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(6, 12).Value = 29
$ws.Cells.Item(7, 12).Value = 134

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(2, 12).Value = 110
$ws.Cells.Item(6, 12).Value = 69
$ws.Cells.Item(7, 12).Value = 289

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Cells.Item(2, 12).Value = 30
$ws.Cells.Item(3, 12).Value = 40
$ws.Cells.Item(7, 12).Value = 101

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(3, 12).Value = 68
$ws.Cells.Item(4, 12).Value = 10
$ws.Cells.Item(7, 12).Value = 192

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(2, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 154

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(2, 12).Value = 70
$ws.Cells.Item(7, 12).Value = 223

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Cells.Item(6, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 62

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Cells.Item(3, 12).Value = 41
$ws.Cells.Item(4, 12).Value = 6
$ws.Cells.Item(7, 12).Value = 113

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Cells.Item(2, 12).Value = 27
$ws.Cells.Item(3, 12).Value = 41
$ws.Cells.Item(7, 12).Value = 130

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(3, 12).Value = 360
$ws.Cells.Item(7, 12).Value = 880

$ws = $wb.Worksheets.Item("Oakland")
$ws.Cells.Item(2, 12).Value = 22
$ws.Cells.Item(7, 12).Value = 49

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(2, 12).Value = 20
$ws.Cells.Item(3, 12).Value = 17

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Cells.Item(2, 12).Value = 42
$ws.Cells.Item(7, 12).Value = 117

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Cells.Item(6, 12).Value = 8
$ws.Cells.Item(7, 12).Value = 32

$ws = $wb.Worksheets.Item("Little Village")
$ws.Cells.Item(3, 12).Value = 117
$ws.Cells.Item(7, 12).Value = 354
